$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$cell = $ws.Range("A3")
$cell.Value = "Chris,`nChambers"
$cell.WrapText = $true
$ws.Rows.Item(3).RowHeight = 29
$ws.Range("A4").Select()
